$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 29.07024133333333
$ws.Range("H2").Value = 87.210724
$ws.Range("I2").Value = 0.4171675701339755
$ws.Range("J2").Value = 0.4171675701339755
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.08181366666667
$ws.Range("N2").Value = 63.245441
$ws.Range("O2").Value = 0.0571606014598545
$ws.Range("P2").Value = 0.0571606014598545
$ws.Range("Q2").Value = 612.8534110343649
$ws.Range("R2").Value = 5515.680699309284
$ws.Range("S2").Value = 0.02384554921840408
$ws.Range("T2").Value = 0.02384554921840407
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 29.07024133333333
$ws.Range("H3").Value = 87.210724
$ws.Range("I3").Value = 0.4171675701339755
$ws.Range("J3").Value = 0.4171675701339755
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8177496571571792
$ws.Range("P3").Value = 0.8177496571571792
$ws.Range("Q3").Value = 8767.589107909211
$ws.Range("R3").Value = 78908.3019711829
$ws.Range("S3").Value = 0.341138637454152
$ws.Range("T3").Value = 0.341138637454152
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 29.07024133333333
$ws.Range("H4").Value = 87.210724
$ws.Range("I4").Value = 0.4171675701339755
$ws.Range("J4").Value = 0.4171675701339755
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1250897413829664
$ws.Range("P4").Value = 0.1250897413829664
$ws.Range("Q4").Value = 1341.162841783586
$ws.Range("R4").Value = 12070.46557605228
$ws.Range("S4").Value = 0.05218338346141948
$ws.Range("T4").Value = 0.05218338346141947
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 28.10308133333334
$ws.Range("H5").Value = 84.30924400000001
$ws.Range("I5").Value = 0.4032885045113541
$ws.Range("J5").Value = 0.4032885045113541
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.08181366666667
$ws.Range("N5").Value = 63.245441
$ws.Range("O5").Value = 0.0571606014598545
$ws.Range("P5").Value = 0.0571606014598545
$ws.Range("Q5").Value = 592.4639241285116
$ws.Range("R5").Value = 5332.175317156604
$ws.Range("S5").Value = 0.02305221347971425
$ws.Range("T5").Value = 0.02305221347971425
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 28.10308133333334
$ws.Range("H6").Value = 84.30924400000001
$ws.Range("I6").Value = 0.4032885045113541
$ws.Range("J6").Value = 0.4032885045113541
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 301.6001486666667
$ws.Range("N6").Value = 904.800446
$ws.Range("O6").Value = 0.8177496571571792
$ws.Range("P6").Value = 0.8177496571571792
$ws.Range("Q6").Value = 8475.893508124758
$ws.Range("R6").Value = 76283.04157312283
$ws.Range("S6").Value = 0.3297890362995913
$ws.Range("T6").Value = 0.3297890362995913
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 28.10308133333334
$ws.Range("H7").Value = 84.30924400000001
$ws.Range("I7").Value = 0.4032885045113541
$ws.Range("J7").Value = 0.4032885045113541
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.13524966666667
$ws.Range("N7").Value = 138.405749
$ws.Range("O7").Value = 0.1250897413829664
$ws.Range("P7").Value = 0.1250897413829664
$ws.Range("Q7").Value = 1296.542673715973
$ws.Range("R7").Value = 11668.88406344376
$ws.Range("S7").Value = 0.05044725473204855
$ws.Range("T7").Value = 0.05044725473204855
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.51148366666667
$ws.Range("H8").Value = 37.534451
$ws.Range("I8").Value = 0.1795439253546705
$ws.Range("J8").Value = 0.1795439253546705
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.08181366666667
$ws.Range("N8").Value = 63.245441
$ws.Range("O8").Value = 0.0571606014598545
$ws.Range("P8").Value = 0.0571606014598545
$ws.Range("Q8").Value = 263.7647673542102
$ws.Range("R8").Value = 2373.882906187891
$ws.Range("S8").Value = 0.01026283876173619
$ws.Range("T8").Value = 0.01026283876173618
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.51148366666667
$ws.Range("H9").Value = 37.534451
$ws.Range("I9").Value = 0.1795439253546705
$ws.Range("J9").Value = 0.1795439253546705
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 301.6001486666667
$ws.Range("N9").Value = 904.800446
$ws.Range("O9").Value = 0.8177496571571792
$ws.Range("P9").Value = 0.8177496571571792
$ws.Range("Q9").Value = 3773.465333907239
$ws.Range("R9").Value = 33961.18800516515
$ws.Range("S9").Value = 0.146821983403436
$ws.Range("T9").Value = 0.146821983403436
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 12.51148366666667
$ws.Range("H10").Value = 37.534451
$ws.Range("I10").Value = 0.1795439253546705
$ws.Range("J10").Value = 0.1795439253546705
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.13524966666667
$ws.Range("N10").Value = 138.405749
$ws.Range("O10").Value = 0.1250897413829664
$ws.Range("P10").Value = 0.1250897413829664
$ws.Range("Q10").Value = 577.2204226620889
$ws.Range("R10").Value = 5194.9838039588
$ws.Range("S10").Value = 0.02245910318949835
$ws.Range("T10").Value = 0.02245910318949835
